$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.216.51"
$ws.Range("E2").Value = "  -2.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.821.88"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.50"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4255"
$ws.Range("E7").Value = "  -2.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3681"
$ws.Range("E8").Value = "  -2.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07238"
$ws.Range("E9").Value = "  -2.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8609"
$ws.Range("E10").Value = "  -2.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.99"
$ws.Range("E11").Value = "  -3.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.825.32"
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.687"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07110"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.311"
$ws.Range("E15").Value = "  -3.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.13"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008869"
$ws.Range("E18").Value = "  -1.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.06"
$ws.Range("E20").Value = "  -2.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.240.90"
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.136"
$ws.Range("E22").Value = "  -2.60%  "
$ws.Range("E23").Value = "  -3.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.048.41"
$ws.Range("E24").Value = "  -1.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.003"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.30"
$ws.Range("E26").Value = "  -2.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.36"
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.124"
$ws.Range("E28").Value = "  +5.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.217"
$ws.Range("E29").Value = "  -3.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.23"
$ws.Range("E30").Value = "  -3.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08883"
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.195"
$ws.Range("E32").Value = "  -3.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7568"
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.432"
$ws.Range("E34").Value = "  -3.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.806"
$ws.Range("E35").Value = "  -6.24%  "
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("E37").Value = "  -2.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01972"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.128"
$ws.Range("E40").Value = "  +1.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.870"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1689"
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5030"
$ws.Range("E43").Value = "  -3.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.621"
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.65"
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.93"
$ws.Range("E46").Value = "  -3.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4739"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06369"
$ws.Range("E49").Value = "  -1.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.660"
$ws.Range("E50").Value = "  -3.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.806"
$ws.Range("E51").Value = "  -3.63%  "
